# Minor changes to instructors' students view
#
# A new subject record was effectively added for student
# "dela Cruz Arvin Sotto" (a Pending course, handled by faculty member
# "Capco Anthony Christian Wee", course CS413-M, for program BSCS /
# college COS, period 2024-2nd semester, Male). This cascades through
# several summary sheets.

$wb = $excel.ActiveWorkbook

# 1. "Number of Enrollees Per Program" - BSCS count 2 -> 3
$wsProgram = $wb.Worksheets.Item("Number of Enrollees Per Program")
$wsProgram.Range("B2").Value = 3

# 2. "Per College" - COS count 3 -> 4
$wsCollege = $wb.Worksheets.Item("Per College")
$wsCollege.Range("B2").Value = 4

# 3. "Per Period" - "2024-2nd semester-COS-BSCS" count 0 -> 1
$wsPeriod = $wb.Worksheets.Item("Per Period")
$wsPeriod.Range("B50").Value = 1

# 4. "Gender" - "2024-2nd semester" Male Count 0 -> 1
$wsGender = $wb.Worksheets.Item("Gender")
$wsGender.Range("B5").Value = 1

# 5. "Faculty - Subjects" - Capco Anthony Christian Wee now handles an
#    extra CS413-M course (duplicate), so courses handled list, pending
#    count, and courses count all grow by one. Widen column B to fit
#    the longer text.
$wsFaculty = $wb.Worksheets.Item("Faculty - Subjects")
# Excel's ColumnWidth property (chars) differs from the raw OOXML
# <col width="..."> units by a per-pixel padding factor; 56.17 is the
# ColumnWidth that serializes back out to width="57".
$wsFaculty.Columns.Item(2).ColumnWidth = 56.17
$wsFaculty.Range("B2").Value = "GEE12D-M, CS413-M, CS413-M, CC413-M, GEE12D-M, GEE11D-M"
$wsFaculty.Range("C2").Value = 2
$wsFaculty.Range("F2").Value = 6

# 6. "Student - Subjects" - "dela Cruz Arvin Sotto" gains a Pending
#    subject, raising Pending and Total Subjects counts by one.
$wsStudent = $wb.Worksheets.Item("Student - Subjects")
$wsStudent.Range("B7").Value = 1
$wsStudent.Range("E7").Value = 1
